# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
# Columns: A=rank index, B=Coin, C=Link, D=Price, E=Volume(1h).
# D/E are stored as text in the sheet, so numeric-looking D values are
# written with a leading apostrophe to keep Excel from coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.603.11'
$ws.Range('E2').Value = '  -2.30%  '
$ws.Range('D3').Value = '''1.587.41'
$ws.Range('E3').Value = '  -2.69%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''211.06'
$ws.Range('E5').Value = '  -2.33%  '
$ws.Range('E6').Value = '  -2.65%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.248'
$ws.Range('E8').Value = '  -2.81%  '
$ws.Range('D9').Value = '''0.0617'
$ws.Range('E9').Value = '  -1.30%  '
$ws.Range('E10').Value = '  -3.18%  '
$ws.Range('D11').Value = '''0.0836'
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('D12').Value = '''1.809.46'
$ws.Range('E12').Value = '  -2.69%  '
$ws.Range('D13').Value = '''1.593.75'
$ws.Range('E13').Value = '  -2.22%  '
$ws.Range('D14').Value = '''4.04'
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('D15').Value = '''0.527'
$ws.Range('E15').Value = '  -3.03%  '
$ws.Range('D16').Value = '''64.78'
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').Value = '''26.607.40'
$ws.Range('E17').Value = '  -2.22%  '
$ws.Range('E18').Value = '  -0.87%  '
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = '''206.99'
$ws.Range('E20').Value = '  -4.47%  '
$ws.Range('E21').Value = '  -2.87%  '
$ws.Range('E22').Value = '  -3.22%  '
$ws.Range('E23').Value = '  -3.33%  '
$ws.Range('D24').Value = '''8.88'
$ws.Range('E24').Value = '  -2.68%  '
$ws.Range('D25').Value = '''147.37'
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E28').Value = '  -3.07%  '
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('D32').Value = '''3.25'
$ws.Range('E32').Value = '  -4.11%  '
$ws.Range('D33').Value = '''0.663'
$ws.Range('E33').Value = '  +22.88%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').Value = '''1.328.11'
$ws.Range('E34').Value = '  +0.79%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '''2.92'
$ws.Range('E35').Value = '  -2.96%  '
$ws.Range('D36').Value = '''1.51'
$ws.Range('E36').Value = '  -3.18%  '
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('E38').Value = '  -1.50%  '
$ws.Range('D39').Value = '''0.826'
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('D41').Value = '''5.41'
$ws.Range('E41').Value = '  +4.18%  '
$ws.Range('D42').Value = '''0.784'
$ws.Range('E42').Value = '  -1.88%  '
$ws.Range('E43').Value = '  -3.67%  '
$ws.Range('D44').Value = '''63.45'
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('D45').Value = '''1.722.29'
$ws.Range('E45').Value = '  -2.54%  '
$ws.Range('D46').Value = '''89.92'
$ws.Range('E46').Value = '  -0.95%  '
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('D48').Value = '''0.829'
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('E49').Value = '  -1.76%  '
$ws.Range('D50').Value = '''0.0979'
$ws.Range('E50').Value = '  +2.40%  '
$ws.Range('D51').Value = '''7.48'
$ws.Range('E51').Value = '  -0.89%  '
